# Commit: "#5: property boat&car done"
#
# This edit:
#  1) Rebuilds the "汽車" (car) sheet (sheet index 3) from a minimal 4-column
#     table into the full 14-column schema shared by the other property
#     sheets (name, capacity, owner, register_date, register_reason,
#     acquire_value, property_category, category, date, legislator_name,
#     legislator_id, source_file, index), filling in the one data row with
#     the car's capacity, owner, registration date/reason.
#  2) Renumbers the running "index" column (col A, and col N where present)
#     on every later sheet (存款/股票/基金受益憑證/債權/事業投資) to stay
#     consistent with the row that was inserted earlier in the global
#     dataset (the car row grew, and an out-of-workbook "boat" dataset
#     shrank by one row upstream of these sheets).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Sheet 3 "汽車" (car) — rebuild header + data row
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

# Stretch the formatting of the existing header/data cells across the new
# columns first (Copy carries over style + number format), then overwrite
# the values. B1/B2 already hold the right style (s=1 header / s=2 data).
$ws3.Range("B1").Copy($ws3.Range("C1:N1"))
$ws3.Range("B2").Copy($ws3.Range("C2:N2"))

# -- header row 1 --
$ws3.Range("B1").Value = "name"
$ws3.Range("C1").Value = "capacity"
$ws3.Range("D1").Value = "owner"
$ws3.Range("E1").Value = "register_date"
$ws3.Range("F1").Value = "register_reason"
$ws3.Range("G1").Value = "acquire_value"
$ws3.Range("H1").Value = "property_category"
$ws3.Range("I1").Value = "category"
$ws3.Range("J1").Value = "date"
$ws3.Range("K1").Value = "legislator_name"
$ws3.Range("L1").Value = "legislator_id"
$ws3.Range("M1").Value = "source_file"
$ws3.Range("N1").Value = "index"

# -- data row 2 --
$ws3.Range("A2").Value = 32
$ws3.Range("B2").Value = "富豪VOLVO"
$ws3.Range("C2").Value = 2435
$ws3.Range("D2").Value = "吳虹"
$ws3.Range("E2").Value = "91年10月07日"
$ws3.Range("F2").Value = "貿賣"
$ws3.Range("G2").Value = "(超過五年）"
$ws3.Range("H2").Value = "land"
$ws3.Range("I2").Value = "normal"
$ws3.Range("J2").Value = "2012-04-25"
$ws3.Range("K2").Value = "陳其邁"
$ws3.Range("L2").Value = 734
$ws3.Range("M2").Value = "tmpa37b1"
$ws3.Range("N2").Value = 32

# ---------------------------------------------------------------------
# 2) Renumber the "index" column on the remaining sheets
# ---------------------------------------------------------------------

# Sheet 4 "存款" (deposit) — column A, rows 2-13: 50..61 -> 47..58
$ws4 = $wb.Worksheets.Item(4)
$vals4 = @(47,48,49,50,51,52,53,54,55,56,57,58)
for ($i = 0; $i -lt $vals4.Length; $i++) {
    $ws4.Cells.Item($i + 2, 1).Value = $vals4[$i]
}

# Sheet 5 "股票" (stock) — columns A & N, rows 2-4: 68..70 -> 65..67
$ws5 = $wb.Worksheets.Item(5)
$vals5 = @(65,66,67)
for ($i = 0; $i -lt $vals5.Length; $i++) {
    $ws5.Cells.Item($i + 2, 1).Value = $vals5[$i]
    $ws5.Cells.Item($i + 2, 14).Value = $vals5[$i]
}

# Sheet 6 "基金受益憑證" (fund) — column A, rows 2-11: 80..89 -> 77..86
$ws6 = $wb.Worksheets.Item(6)
$vals6 = @(77,78,79,80,81,82,83,84,85,86)
for ($i = 0; $i -lt $vals6.Length; $i++) {
    $ws6.Cells.Item($i + 2, 1).Value = $vals6[$i]
}

# Sheet 7 "債權" (debt) — column A, rows 2-7: 110..115 -> 107..112
$ws7 = $wb.Worksheets.Item(7)
$vals7 = @(107,108,109,110,111,112)
for ($i = 0; $i -lt $vals7.Length; $i++) {
    $ws7.Cells.Item($i + 2, 1).Value = $vals7[$i]
}

# Sheet 8 "事業投資" (business investment) — column A, rows 2-3: 120,121 -> 117,118
$ws8 = $wb.Worksheets.Item(8)
$vals8 = @(117,118)
for ($i = 0; $i -lt $vals8.Length; $i++) {
    $ws8.Cells.Item($i + 2, 1).Value = $vals8[$i]
}
